$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '42.817.22'
$c.Style = $s
$ws.Range("E2").Value = '  +0.86%  '
$c = $ws.Range("D3")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.526.16'
$c.Style = $s
$ws.Range("E3").Value = '  +0.26%  '
$c = $ws.Range("D4")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = $s
$ws.Range("E4").Value = '  -0.14%  '
$c = $ws.Range("D5")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '316.46'
$c.Style = $s
$ws.Range("E5").Value = '  +4.43%  '
$c = $ws.Range("D6")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '94.77'
$c.Style = $s
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("E9").Value = '  -0.93%  '
$c = $ws.Range("D10")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '35.92'
$c.Style = $s
$ws.Range("E10").Value = '  -0.96%  '
$c = $ws.Range("D11")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0810'
$c.Style = $s
$ws.Range("E11").Value = '  +0.43%  '
$c = $ws.Range("D12")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '7.58'
$c.Style = $s
$ws.Range("E12").Value = '  -0.61%  '
$ws.Range("E13").Value = '  -1.68%  '
$c = $ws.Range("D14")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.914.84'
$c.Style = $s
$ws.Range("E14").Value = '  +0.21%  '
$c = $ws.Range("D15")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.517.96'
$c.Style = $s
$ws.Range("E15").Value = '  +0.09%  '
$c = $ws.Range("D16")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '15.22'
$c.Style = $s
$ws.Range("E16").Value = '  +1.83%  '
$ws.Range("E17").Value = '  -1.36%  '
$c = $ws.Range("D18")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '42.904.28'
$c.Style = $s
$ws.Range("E18").Value = '  +0.97%  '
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("E20").Value = '  +3.78%  '
$ws.Range("E21").Value = '  -0.51%  '
$c = $ws.Range("D22")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '69.90'
$c.Style = $s
$ws.Range("E22").Value = '  -1.45%  '
$c = $ws.Range("D23")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '250.92'
$c.Style = $s
$ws.Range("E23").Value = '  +0.33%  '
$c = $ws.Range("D24")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = $s
$ws.Range("E24").Value = '  +2.28%  '
$ws.Range("E25").Value = '  +0.35%  '
$c = $ws.Range("D26")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '26.77'
$c.Style = $s
$ws.Range("E26").Value = '  -0.71%  '
$ws.Range("E27").Value = '  -0.08%  '
$c = $ws.Range("D28")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.Style = $s
$ws.Range("E28").Value = '  +4.04%  '
$c = $ws.Range("D29")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '40.00'
$c.Style = $s
$ws.Range("E29").Value = '  +4.85%  '
$ws.Range("E30").Value = '  +0.17%  '
$c = $ws.Range("D31")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.00'
$c.Style = $s
$ws.Range("E31").Value = '  +1.34%  '
$c = $ws.Range("D32")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '155.38'
$c.Style = $s
$ws.Range("E32").Value = '  +0.50%  '
$ws.Range("E33").Value = '  +3.20%  '
$c = $ws.Range("D34")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '19.09'
$c.Style = $s
$ws.Range("E34").Value = '  +2.87%  '
$ws.Range("E35").Value = '  -0.48%  '
$c = $ws.Range("D36")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0788'
$c.Style = $s
$ws.Range("E36").Value = '  +0.44%  '
$c = $ws.Range("D37")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.60'
$c.Style = $s
$ws.Range("E37").Value = '  -0.77%  '
$ws.Range("E38").Value = '  -2.21%  '
$ws.Range("E39").Value = '  +0.04%  '
$c = $ws.Range("D40")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '23.64'
$c.Style = $s
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("E41").Value = '  +14.18%  '
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("E43").Value = '  +0.20%  '
$c = $ws.Range("D44")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '3.77'
$c.Style = $s
$ws.Range("E44").Value = '  -1.05%  '
$ws.Range("E45").Value = '  -2.04%  '
$c = $ws.Range("D46")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.018.19'
$c.Style = $s
$ws.Range("E46").Value = '  -0.33%  '
$c = $ws.Range("D47")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '85.53'
$c.Style = $s
$ws.Range("E47").Value = '  +1.19%  '
$c = $ws.Range("D48")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.78'
$c.Style = $s
$ws.Range("E48").Value = '  -1.55%  '
$c = $ws.Range("D49")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.769.50'
$c.Style = $s
$ws.Range("E49").Value = '  +0.11%  '
$c = $ws.Range("D50")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '73.36'
$c.Style = $s
$ws.Range("E50").Value = '  +2.32%  '
$c = $ws.Range("D51")
$s = $c.Style
$c.NumberFormat = "@"
$c.Value = '102.40'
$c.Style = $s
$ws.Range("E51").Value = '  +1.02%  '
